# Add data for 2021-12-21 (Dec 13 cutoff, week of Dec 21 observations)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet and update the "through" label to reflect new cutoff date
$ws.Name = "Through 2021-12-13"
$ws.Range("B1").Value = "December 2021 (through December 13)"

# West Town (row 2) - new observation
$ws.Range("BV2").Value = 1

# North Lawndale (row 4)
$ws.Range("B4").Value = 4
$ws.Range("N4").Value = 9

# Garfield Park (row 6)
$ws.Range("N6").Value = 9
$ws.Range("AX6").Value = 5
$ws.Range("BJ6").Value = 2

# Austin (row 7)
$ws.Range("AX7").Value = 3

# Chatham (row 8)
$ws.Range("BJ8").Value = 2

# Douglas (row 10)
$ws.Range("Z10").Value = 1

# Chicago Lawn (row 14)
$ws.Range("AL14").Value = 1

# Grand Boulevard (row 18)
$ws.Range("B18").Value = 4

# Wicker Park (row 21)
$ws.Range("Z21").Value = 1

# South Shore (row 24)
$ws.Range("N24").Value = 3
$ws.Range("BJ24").Value = 3

# Avalon Park (row 29)
$ws.Range("Z29").Value = 1

# New City (row 32)
$ws.Range("Z32").Value = 2

# Chinatown (row 41)
$ws.Range("B41").Value = 4
$ws.Range("D41").Value = 2

# Rogers Park (row 54)
$ws.Range("N54").Value = 1

# Lincoln Park (row 61)
$ws.Range("B61").Value = 2
$ws.Range("N61").Value = 1

# Archer Heights (row 65)
$ws.Range("B65").Value = 1

# Belmont Cragin (row 67)
$ws.Range("N67").Value = 1

# North Center (row 88)
$ws.Range("AX88").Value = 2
